# Add a new QoQ forecast column (BB) covering one more quarter, and a new
# trailing row (83) for the additional period, per the commit:
# "Included EQUIPMENT eval, updated folder structure, allowed for multiple
# archive excels".
#
# Structure of the sheet:
#  - Row 1 (header, style "1"): per-vintage forecast-origin dates in B1:BA1.
#    BB1 is the new vintage date, one quarter after BA1.
#  - Column A (style "1"): target-period dates, one per data row.
#  - Data rows 2-82: forecast values per vintage column. BB matches BA for
#    rows 2-70 (unchanged forecasts) and gets new, revised values for rows
#    71-82 (the new vintage revises the recent quarters).
#  - Row 83 is brand new: one quarter after row 82, with only column A
#    (date) and the new BB column populated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 header: new vintage date in BB1, matching A-column/header style ---
$ws.Range("BA1").Copy()
$ws.Range("BB1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("BB1").Value = 45986

# --- Rows 2-70: BB repeats the same value already in BA (no style) ---
for ($r = 2; $r -le 70; $r++) {
    $ws.Cells.Item($r, 54).Value = $ws.Cells.Item($r, 53).Value()
}

# --- Rows 71-82: BB gets the revised forecast values from the new vintage ---
$revised = @{
    71 = -0.1118837721692358
    72 = 0.3266766184601977
    73 = 0.325608361860148
    74 = 0.2270536959888376
    75 = 0.2584586613899786
    76 = 0.267110162551939
    77 = 0.2584457304690463
    78 = 0.25955667089935
    79 = 0.2608396828010808
    80 = 0.2602238633310655
    81 = 0.2601767679767295
    82 = 0.2603078675382955
}
foreach ($r in $revised.Keys) {
    $ws.Cells.Item($r, 54).Value = $revised[$r]
}

# --- Row 83: brand new trailing row, one quarter after row 82 ---
$ws.Range("A82").Copy()
$ws.Range("A83").PasteSpecial(-4122)   # xlPasteFormats (match date-column style)
$ws.Range("A83").Value = 46934
$ws.Range("BB83").Value = 0.2602752457138798
